$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily work-tracking entry in row 23 (Date 13.12)
$ws.Range("A23").Value = 13.12

# Time IN / Time OUT - copy the existing time-format from the row above so
# the cells keep using the workbook's existing time style instead of a new one
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B23").Value = 0.45833333333333331   # 11:00 AM

$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C23").Value = 0.64583333333333337   # 3:30 PM

$excel.CutCopyMode = $false

# Sum. Time and Activities text for the new entry
$ws.Range("E23").Value = "4hr 30min"
$ws.Range("F23").Value = "puttig html layout more and made combinatorics differently and random value make sort"

# Move the active cell selection to F23, matching the edited sheet view
$ws.Range("F23").Select()
